# Apply the design/theme change: switch the presentation's colour theme
# from "Integral" (Red Violet) to the built-in "Office Theme" colours,
# and update the table on slide 5 to use the corresponding built-in
# "Office" table style.

$p = $ppt.ActivePresentation

# --- 1) Re-colour the presentation theme (Design > Colors) ------------
# The font scheme and format scheme (fills/lines/effects) are already
# identical between the two themes - only the 12 theme colours differ,
# so recolouring the active theme reproduces the target theme exactly.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Order of ThemeColorScheme items: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink - matches the RGB() values (R + G*256 + B*65536) below,
# taken from the "Office" theme palette.
$officeRgb = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeRgb[$i - 1]
}

# --- 2) Re-style the table on slide 5 (Table Design > Table Styles) ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{78275AF7-A210-4F55-A953-33676744B367}")
}
